$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$data = @(
    @(2, "VAIBHAVGBL", 11.2793, 20.0938, 27.2582),
    @(3, "KELLTONTEC", 11.1389, 8.167199999999999, -6.7868),
    @(4, "ASALCBR", 10.4199, 10.6262, 24.6814),
    @(5, "SALASAR", 10.1523, 15.4255, 22.3224),
    @(6, "SAGILITY", 9.9862, 17.7399, 31.1533),
    @(7, "UNIPARTS", 8.911300000000001, 11.2136, 27.5883),
    @(8, "MARINE", 7.3765, 3.8718, 16.5561),
    @(9, "CHENNPETRO", 7.1508, 12.3325, 14.4893),
    @(10, "MCLOUD", 7.0251, 5.9993, -22.4516),
    @(11, "INDIACEM", 5.6343, 5.9185, 7.6562),
    @(12, "SHREEJISPG", 5.5842, 9.837300000000001, 10.9267),
    @(13, "PDSL", 5.4931, 8.8353, 14.9929),
    @(14, "VENKEYS", 5.4574, 6.087, 3.7795),
    @(15, "CREDITACC", 5.4499, 4.0562, 9.385199999999999),
    @(16, "NETWEB", 5.4177, 11.2911, 13.3836),
    @(17, "BLSE", 5.2743, 4.2419, -1.9432),
    @(18, "FIVESTAR", 5.2723, 18.6127, 18.7011),
    @(19, "BHEL", 5.1184, 11.6957, 8.1778),
    @(20, "POLICYBZR", 5.0714, 7.4213, 6.3925),
    @(21, "SOLEX", 4.9978, 6.5692, "N/A"),
    @(22, "V2RETAIL", 4.9975, 3.4793, 11.7421),
    @(23, "SKMEGGPROD", 4.9909, 12.0154, 29.9408),
    @(24, "BAJAJHCARE", 4.8905, 5.4418, -0.8851),
    @(25, "MTARTECH", 4.7599, 8.930999999999999, 33.1425),
    @(26, "VIMTALABS", 4.7315, 4.9107, -0.2123),
    @(27, "EUROPRATIK", 4.6535, 10.5452, 27.5556),
    @(28, "RELTD", 4.4461, 11.3145, -0.3467),
    @(29, "MRPL", 4.3533, 14.4863, 25.2805),
    @(30, "DEEDEV", 4.2721, -2.6649, -3.4677),
    @(31, "CENTRUM", 4.1515, 2.7504, 2.079),
    @(32, "VSTIND", 4.0418, 4.5472, 4.082),
    @(33, "RAMASTEEL", 4.0161, 3.9117, 5.499),
    @(34, "SAPPHIRE", 3.9629, 5.7961, 3.1314),
    @(35, "IIFL", 3.8623, 10.9024, 20.2036),
    @(36, "ABREL", 3.8178, 12.1731, 11.7083),
    @(37, "HIRECT", 3.8132, 11.2414, 9.533200000000001),
    @(38, "BLS", 3.7397, 0.6756, -0.5895),
    @(39, "SKYGOLD", 3.7048, -0.8479, 37.7018),
    @(40, "SHANTIGOLD", 3.6177, 10.9107, 3.5111),
    @(41, "PFOCUS", 3.5616, 0.8404, 2.3019),
    @(42, "OIL", 3.5196, 3.7664, 5.1957),
    @(43, "PSPPROJECT", 3.5033, 17.6634, 24.1324),
    @(44, "REDTAPE", 3.4892, 3.3888, -3.5311),
    @(45, "CENTUM", 3.402, 4.2653, -1.2627),
    @(46, "SHRINGARMS", 3.3845, 4.5935, 24.7977),
    @(47, "BGRENERGY", 3.2276, -6.0804, 74.84910000000001),
    @(48, "KMEW", 3.2076, 5.3329, -0.2666),
    @(49, "IVALUE", 3.1867, 6.4398, -0.7706),
    @(50, "ALICON", 3.1655, 9.2813, 14.7096),
    @(51, "RSYSTEMS", 3.1129, 4.2704, 6.6118),
    @(52, "BEML", 3.0957, 0.178, 6.5216),
    @(53, "POWERINDIA", 3.0946, 7.4909, 0.1222),
    @(54, "ASHOKA", 3.0282, 4.5515, 7.1831),
    @(55, "CARYSIL", 3.0237, 2.5061, 11.425),
    @(56, "BPCL", 2.9733, 8.4733, 5.5351),
    @(57, "NBCC", 2.9682, 6.2226, 10.7957),
    @(58, "GMMPFAUDLR", 2.945, 7.4065, 19.7414),
    @(59, "GANESHCP", 2.9391, 2.4132, 1.9642),
    @(60, "NEULANDLAB", 2.869, -1.4918, 8.5108),
    @(61, "DBCORP", 2.7678, 5.4075, 1.3556),
    @(62, "JKTYRE", 2.7467, 5.7895, 21.851),
    @(63, "SPANDANA", 2.6709, 4.1658, 2.7815),
    @(64, "AHLUCONT", 2.6425, 1.4547, -5.7426),
    @(65, "OBEROIRLTY", 2.638, 3.4531, 11.1561),
    @(66, "GREAVESCOT", 2.5953, 15.1557, 11.3996),
    @(67, "REFEX", 2.5106, -0.055, 1.9212),
    @(68, "SDBL", 2.4986, 1.0448, 6.6313),
    @(69, "JSFB", 2.4481, 2.3102, -4.631),
    @(70, "JKLAKSHMI", 2.4379, 4.424, 1.4462),
    @(71, "DBL", 2.4138, 3.4966, 4.5695),
    @(72, "INOXGREEN", 2.3804, 10.3625, 33.6632),
    @(73, "VOLTAMP", 2.3791, 2.2586, 1.9971),
    @(74, "MAMATA", 2.3783, 1.7903, 0.9943),
    @(75, "MFSL", 2.3703, 2.4242, -1.3389),
    @(76, "CAMS", 2.3291, 1.6778, 4.9673)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Applied updates to Top Gainers rows 2-76"
